# Update portfolio data (Web UI sync at 2025-09-24 02:08) across the three
# sheets: 大智投资组合 (sheet1), 大成投资组合 (sheet2), 我的投资组合 (sheet3).
#
# Helper: write a value as genuine TEXT (not auto-coerced to a number),
# matching the source workbook's convention of storing stock codes /
# timestamps as strings. We briefly force a text number-format so the
# COM layer doesn't "smart type" a digit-only string into a number, then
# reset the style back to Normal so we don't leave a stray cell style
# behind (the source cells carry no explicit style).
function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$NEWTS = "202509241007"

# ---------------------------------------------------------------------
# Sheet 1: 大智投资组合
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("大智投资组合")

$ws1.Range("D2").Value = 5.06
Set-TextCell $ws1 "E2" $NEWTS

$ws1.Range("D3").Value = 9.61
Set-TextCell $ws1 "E3" $NEWTS

$ws1.Range("D4").Value = 8.92
Set-TextCell $ws1 "E4" $NEWTS

$ws1.Range("D5").Value = 5.26
Set-TextCell $ws1 "E5" $NEWTS

$ws1.Range("D6").Value = 5.13
Set-TextCell $ws1 "E6" $NEWTS

$ws1.Range("D7").Value = 2.11
Set-TextCell $ws1 "E7" $NEWTS

$ws1.Range("D8").Value = 1.92
Set-TextCell $ws1 "E8" $NEWTS

Set-TextCell $ws1 "E9" $NEWTS

# ---------------------------------------------------------------------
# Sheet 2: 大成投资组合
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("大成投资组合")

Set-TextCell $ws2 "E2" $NEWTS

$ws2.Range("D3").Value = 5.44
Set-TextCell $ws2 "E3" $NEWTS

$ws2.Range("D4").Value = 4.7
Set-TextCell $ws2 "E4" $NEWTS

$ws2.Range("D5").Value = 6.14
Set-TextCell $ws2 "E5" $NEWTS

$ws2.Range("D6").Value = 5.23
Set-TextCell $ws2 "E6" $NEWTS

$ws2.Range("D7").Value = 0.97
Set-TextCell $ws2 "E7" $NEWTS

$ws2.Range("D8").Value = 5.89
Set-TextCell $ws2 "E8" $NEWTS

# Row 12 (605580 / 恒盛能源) drops out of the portfolio entirely; deleting
# it shifts the old row 13 (601598 / 中国外运) up into row 12.
$ws2.Rows.Item(12).Delete()

# The remaining rows 9-12 are re-ranked to match the new allocation order,
# so rewrite them directly to the target state.
Set-TextCell $ws2 "B9" "601598"
Set-TextCell $ws2 "C9" "中国外运"
$ws2.Range("D9").Value = 32.01
Set-TextCell $ws2 "E9" $NEWTS

Set-TextCell $ws2 "B10" "601878"
Set-TextCell $ws2 "C10" "浙商证券"
$ws2.Range("D10").Value = 4.85
Set-TextCell $ws2 "E10" $NEWTS

Set-TextCell $ws2 "B11" "603119"
Set-TextCell $ws2 "C11" "浙江荣泰"
$ws2.Range("D11").Value = 0.03
Set-TextCell $ws2 "E11" $NEWTS

Set-TextCell $ws2 "B12" "HK01810"
Set-TextCell $ws2 "C12" "小米集团-W"
$ws2.Range("D12").Value = 1
Set-TextCell $ws2 "E12" $NEWTS

# ---------------------------------------------------------------------
# Sheet 3: 我的投资组合
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("我的投资组合")

$ws3.Range("F2").Value = 5.06
Set-TextCell $ws3 "G2" $NEWTS

Set-TextCell $ws3 "G3" $NEWTS

$ws3.Range("F4").Value = 5.08
Set-TextCell $ws3 "G4" $NEWTS

$ws3.Range("F5").Value = 9.28
Set-TextCell $ws3 "G5" $NEWTS

$ws3.Range("F6").Value = 6.94
Set-TextCell $ws3 "G6" $NEWTS

Set-TextCell $ws3 "G7" $NEWTS

$ws3.Range("F8").Value = 5.38
Set-TextCell $ws3 "G8" $NEWTS

$ws3.Range("F9").Value = 3.21
Set-TextCell $ws3 "G9" $NEWTS

$ws3.Range("F10").Value = 0.98
Set-TextCell $ws3 "G10" $NEWTS

$ws3.Range("F11").Value = 5.06
Set-TextCell $ws3 "G11" $NEWTS

Set-TextCell $ws3 "G12" $NEWTS

Set-TextCell $ws3 "G13" $NEWTS

Write-Host "portfolio data refreshed"
